$d = $word.ActiveDocument

function Get-ParaIndexByText {
    param(
        [string]$pattern
    )
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

# 1) Insert a new bullet (ilvl 0) "El usuario podrá cambiar la fecha de la habitación"
#    right before the "Usuario administrador:" paragraph.
$idx = Get-ParaIndexByText "^Usuario administrador:"
if ($idx -lt 0) {
    throw "Could not find 'Usuario administrador:' paragraph"
}
$target = $d.Paragraphs($idx)
$r = $target.Range
$r.Collapse(1)
$r.InsertParagraphBefore()
$newPara = $d.Paragraphs($idx)
$newPara.Range.InsertBefore("El usuario podrá cambiar la fecha de la habitación")

# 2) Insert three new bullets (ilvl 1) before "Puede añadir nuevas habitaciones":
#    "Ver usuario registrados", "Modificar disponibilidad", "Aceptar o cancelar reserva"
$idx2 = Get-ParaIndexByText "^Puede añadir nuevas habitaciones"
if ($idx2 -lt 0) {
    throw "Could not find 'Puede añadir nuevas habitaciones' paragraph"
}

$newTexts = @("Ver usuario registrados", "Modificar disponibilidad", "Aceptar o cancelar reserva")
foreach ($t in $newTexts) {
    $target2 = $d.Paragraphs($idx2)
    $r2 = $target2.Range
    $r2.Collapse(1)
    $r2.InsertParagraphBefore()
}
for ($k = 0; $k -lt $newTexts.Length; $k++) {
    $p = $d.Paragraphs($idx2 + $k)
    $p.Range.InsertBefore($newTexts[$k])
}

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
